$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their original text formatting so that
# numeric-looking strings (e.g. "299.92", "0.100") are not coerced into numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.800.44'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.294.62'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '299.92'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.58'
$ws.Range('E6').Value = '  -4.72%  '
$ws.Range('E7').Value = '  -1.35%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -3.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.28'
$ws.Range('E10').Value = '  -5.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0797'
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.10'
$ws.Range('E12').Value = '  -5.07%  '
$ws.Range('E13').Value = '  +2.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '16.76'
$ws.Range('E14').Value = '  +7.31%  '
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.651.86'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.287.20'
$ws.Range('E17').Value = '  -3.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.804'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.722.12'
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0897'
$ws.Range('E20').Value = '  -1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.50'
$ws.Range('E21').Value = '  -2.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.13'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '235.43'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('E25').Value = '  +0.85%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -3.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.26'
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('E29').Value = '  +1.58%  '
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.69'
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.07'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  +3.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.91'
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.38'
$ws.Range('E36').Value = '  -1.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.70'
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0689'
$ws.Range('E38').Value = '  -2.56%  '
$ws.Range('E39').Value = '  -3.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.100'
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('E41').Value = '  -4.59%  '
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.33'
$ws.Range('E43').Value = '  -2.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.986.95'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('E45').Value = '  -2.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.77'
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.46'
$ws.Range('E47').Value = '  -6.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.83'
$ws.Range('E48').Value = '  -3.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.520.65'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.63'
$ws.Range('E50').Value = '  -5.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.54'
$ws.Range('E51').Value = '  -7.73%  '
